{"js": "// Replace the 100 arithmetic-problem texts in the practice-sheet table\n// with their updated values (cell formatting / runs are left untouched;\n// only the w:t text content of each matching run changes).\nconst replacements = [\n  [\"18+9=\", \"12+52=\"],\n  [\"46-33=\", \"28+46=\"],\n  [\"87-87=\", \"44-7=\"],\n  [\"18-11=\", \"27+29=\"],\n  [\"49+16=\", \"48-34=\"],\n  [\"43+54=\", \"22-11=\"],\n  [\"76-52=\", \"32+22=\"],\n  [\"91-48=\", \"32+62=\"],\n  [\"35-25=\", \"89-12=\"],\n  [\"69-36=\", \"60-59=\"],\n  [\"27+34=\", \"38+40=\"],\n  [\"46+53=\", \"66-54=\"],\n  [\"22+15=\", \"59+2=\"],\n  [\"6+41=\", \"30+10=\"],\n  [\"58-48=\", \"53-1=\"],\n  [\"45-10=\", \"88-11=\"],\n  [\"36+2=\", \"45+0=\"],\n  [\"16+81=\", \"64-1=\"],\n  [\"20+16=\", \"58-10=\"],\n  [\"76-45=\", \"18+7=\"],\n  [\"64-14=\", \"9+60=\"],\n  [\"8+38=\", \"18+76=\"],\n  [\"11+22=\", \"43+49=\"],\n  [\"77-50=\", \"12+4=\"],\n  [\"2+24=\", \"66+1=\"],\n  [\"50+48=\", \"19+60=\"],\n  [\"96-80=\", \"23-12=\"],\n  [\"95-88=\", \"6+74=\"],\n  [\"9+39=\", \"69-34=\"],\n  [\"22-2=\", \"63-63=\"],\n  [\"68-24=\", \"87-57=\"],\n  [\"40+4=\", \"27+23=\"],\n  [\"95-46=\", \"69-33=\"],\n  [\"32+43=\", \"49-27=\"],\n  [\"19+8=\", \"12-9=\"],\n  [\"5+85=\", \"54+25=\"],\n  [\"73-0=\", \"50-0=\"],\n  [\"68+11=\", \"23+16=\"],\n  [\"95-74=\", \"93-85=\"],\n  [\"14+13=\", \"27-3=\"],\n  [\"56-19=\", \"39+18=\"],\n  [\"89+9=\", \"28-18=\"],\n  [\"34+51=\", \"7+55=\"],\n  [\"95-11=\", \"42+19=\"],\n  [\"25+27=\", \"84-66=\"],\n  [\"83-47=\", \"59-37=\"],\n  [\"27-20=\", \"83-39=\"],\n  [\"0+54=\", \"98-9=\"],\n  [\"51-49=\", \"75-17=\"],\n  [\"81+3=\", \"15-8=\"],\n  [\"90-5=\", \"43-15=\"],\n  [\"20+42=\", \"16+41=\"],\n  [\"87-18=\", \"47-18=\"],\n  [\"69-44=\", \"47+9=\"],\n  [\"76-64=\", \"41+18=\"],\n  [\"85+5=\", \"3+47=\"],\n  [\"25-24=\", \"70-64=\"],\n  [\"62-20=\", \"12-7=\"],\n  [\"96-50=\", \"15+74=\"],\n  [\"97-28=\", \"41+52=\"],\n  [\"27-26=\", \"36+58=\"],\n  [\"10+69=\", \"35-32=\"],\n  [\"31+67=\", \"31+37=\"],\n  [\"1+13=\", \"92-44=\"],\n  [\"79-19=\", \"89-35=\"],\n  [\"80-79=\", \"36-13=\"],\n  [\"2+13=\", \"58-1=\"],\n  [\"3+54=\", \"45+13=\"],\n  [\"9+10=\", \"7+35=\"],\n  [\"2+68=\", \"93-9=\"],\n  [\"84-37=\", \"64+17=\"],\n  [\"59-21=\", \"65-16=\"],\n  [\"52-9=\", \"47+23=\"],\n  [\"10+88=\", \"73-69=\"],\n  [\"54-14=\", \"73+15=\"],\n  [\"48-40=\", \"45-21=\"],\n  [\"78+11=\", \"15+38=\"],\n  [\"1+74=\", \"40+48=\"],\n  [\"19-15=\", \"98-71=\"],\n  [\"97-47=\", \"65-43=\"],\n  [\"41+54=\", \"37+43=\"],\n  [\"30-23=\", \"31+17=\"],\n  [\"97-68=\", \"32+66=\"],\n  [\"5+43=\", \"8+78=\"],\n  [\"74-49=\", \"89-76=\"],\n  [\"40+18=\", \"70-49=\"],\n  [\"78-49=\", \"92-60=\"],\n  [\"67+23=\", \"33+2=\"],\n  [\"5+54=\", \"7+42=\"],\n  [\"50+28=\", \"30+36=\"],\n  [\"35+37=\", \"89-36=\"],\n  [\"8+10=\", \"12+38=\"],\n  [\"44-9=\", \"11-0=\"],\n  [\"65+28=\", \"6+5=\"],\n  [\"47-33=\", \"87-11=\"],\n  [\"98-83=\", \"55-19=\"],\n  [\"0+56=\", \"95-45=\"],\n  [\"1+6=\", \"76-74=\"],\n  [\"14+77=\", \"25+30=\"],\n  [\"49+14=\", \"41+29=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem texts in the practice-sheet table\n# with their updated values using Word's Find & Replace (Find.Execute).\n# Cell formatting / runs are left untouched; only the matched text changes.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"18+9=\", \"12+52=\")\n    ,@(\"46-33=\", \"28+46=\")\n    ,@(\"87-87=\", \"44-7=\")\n    ,@(\"18-11=\", \"27+29=\")\n    ,@(\"49+16=\", \"48-34=\")\n    ,@(\"43+54=\", \"22-11=\")\n    ,@(\"76-52=\", \"32+22=\")\n    ,@(\"91-48=\", \"32+62=\")\n    ,@(\"35-25=\", \"89-12=\")\n    ,@(\"69-36=\", \"60-59=\")\n    ,@(\"27+34=\", \"38+40=\")\n    ,@(\"46+53=\", \"66-54=\")\n    ,@(\"22+15=\", \"59+2=\")\n    ,@(\"6+41=\", \"30+10=\")\n    ,@(\"58-48=\", \"53-1=\")\n    ,@(\"45-10=\", \"88-11=\")\n    ,@(\"36+2=\", \"45+0=\")\n    ,@(\"16+81=\", \"64-1=\")\n    ,@(\"20+16=\", \"58-10=\")\n    ,@(\"76-45=\", \"18+7=\")\n    ,@(\"64-14=\", \"9+60=\")\n    ,@(\"8+38=\", \"18+76=\")\n    ,@(\"11+22=\", \"43+49=\")\n    ,@(\"77-50=\", \"12+4=\")\n    ,@(\"2+24=\", \"66+1=\")\n    ,@(\"50+48=\", \"19+60=\")\n    ,@(\"96-80=\", \"23-12=\")\n    ,@(\"95-88=\", \"6+74=\")\n    ,@(\"9+39=\", \"69-34=\")\n    ,@(\"22-2=\", \"63-63=\")\n    ,@(\"68-24=\", \"87-57=\")\n    ,@(\"40+4=\", \"27+23=\")\n    ,@(\"95-46=\", \"69-33=\")\n    ,@(\"32+43=\", \"49-27=\")\n    ,@(\"19+8=\", \"12-9=\")\n    ,@(\"5+85=\", \"54+25=\")\n    ,@(\"73-0=\", \"50-0=\")\n    ,@(\"68+11=\", \"23+16=\")\n    ,@(\"95-74=\", \"93-85=\")\n    ,@(\"14+13=\", \"27-3=\")\n    ,@(\"56-19=\", \"39+18=\")\n    ,@(\"89+9=\", \"28-18=\")\n    ,@(\"34+51=\", \"7+55=\")\n    ,@(\"95-11=\", \"42+19=\")\n    ,@(\"25+27=\", \"84-66=\")\n    ,@(\"83-47=\", \"59-37=\")\n    ,@(\"27-20=\", \"83-39=\")\n    ,@(\"0+54=\", \"98-9=\")\n    ,@(\"51-49=\", \"75-17=\")\n    ,@(\"81+3=\", \"15-8=\")\n    ,@(\"90-5=\", \"43-15=\")\n    ,@(\"20+42=\", \"16+41=\")\n    ,@(\"87-18=\", \"47-18=\")\n    ,@(\"69-44=\", \"47+9=\")\n    ,@(\"76-64=\", \"41+18=\")\n    ,@(\"85+5=\", \"3+47=\")\n    ,@(\"25-24=\", \"70-64=\")\n    ,@(\"62-20=\", \"12-7=\")\n    ,@(\"96-50=\", \"15+74=\")\n    ,@(\"97-28=\", \"41+52=\")\n    ,@(\"27-26=\", \"36+58=\")\n    ,@(\"10+69=\", \"35-32=\")\n    ,@(\"31+67=\", \"31+37=\")\n    ,@(\"1+13=\", \"92-44=\")\n    ,@(\"79-19=\", \"89-35=\")\n    ,@(\"80-79=\", \"36-13=\")\n    ,@(\"2+13=\", \"58-1=\")\n    ,@(\"3+54=\", \"45+13=\")\n    ,@(\"9+10=\", \"7+35=\")\n    ,@(\"2+68=\", \"93-9=\")\n    ,@(\"84-37=\", \"64+17=\")\n    ,@(\"59-21=\", \"65-16=\")\n    ,@(\"52-9=\", \"47+23=\")\n    ,@(\"10+88=\", \"73-69=\")\n    ,@(\"54-14=\", \"73+15=\")\n    ,@(\"48-40=\", \"45-21=\")\n    ,@(\"78+11=\", \"15+38=\")\n    ,@(\"1+74=\", \"40+48=\")\n    ,@(\"19-15=\", \"98-71=\")\n    ,@(\"97-47=\", \"65-43=\")\n    ,@(\"41+54=\", \"37+43=\")\n    ,@(\"30-23=\", \"31+17=\")\n    ,@(\"97-68=\", \"32+66=\")\n    ,@(\"5+43=\", \"8+78=\")\n    ,@(\"74-49=\", \"89-76=\")\n    ,@(\"40+18=\", \"70-49=\")\n    ,@(\"78-49=\", \"92-60=\")\n    ,@(\"67+23=\", \"33+2=\")\n    ,@(\"5+54=\", \"7+42=\")\n    ,@(\"50+28=\", \"30+36=\")\n    ,@(\"35+37=\", \"89-36=\")\n    ,@(\"8+10=\", \"12+38=\")\n    ,@(\"44-9=\", \"11-0=\")\n    ,@(\"65+28=\", \"6+5=\")\n    ,@(\"47-33=\", \"87-11=\")\n    ,@(\"98-83=\", \"55-19=\")\n    ,@(\"0+56=\", \"95-45=\")\n    ,@(\"1+6=\", \"76-74=\")\n    ,@(\"14+77=\", \"25+30=\")\n    ,@(\"49+14=\", \"41+29=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
